$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Construct', ['Token Artifact Creature — Construct', 'Trample', '6/12'])"
$ws.Range("A3").Value = "('Dragon', ['Token Creature — Dragon', 'Flying', '5/5'])"
$ws.Range("A4").Value = "('Elephant', ['Token Creature — Elephant', '3/3'])"
$ws.Range("A5").Value = "('Ogre', ['Token Creature — Ogre', '3/3'])"
$ws.Range("A6").Value = "('Plant', ['Token Creature — Plant', '0/1'])"
$ws.Range("A7").Value = "('Soldier Ally', ['Token Creature — Soldier Ally', '1/1'])"

$ws.Range("A8:A21").EntireRow.Delete()
